# Remove the BOM row for designator "N1" (treceit / Housing), which sat at
# worksheet row 17. Deleting the entire row shifts every subsequent row
# (R1,R2 ... Y1) up by one, matching the target layout (rows 17-29) and
# shrinking the used range from A1:E30 to A1:E29.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Delete() | Out-Null

# Match the author's updated selection, which now points at the row that
# took over position 17 (previously row 18, "R1, R2").
$ws.Range("A17").Select() | Out-Null
